$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the manual highlight formatting that had been applied to rows
#    79:83 (this is the range that was selected: A79:XFD83). The highlight
#    used a one-off fill/style set (fillId 7 / xf 8,9,10) that only these
#    rows referenced. Clearing it and re-matching the normal per-column
#    formatting used by the rest of the "PL Admin Tools" task rows (e.g.
#    rows 41:45, which share the same Task Name) restores the cells to
#    their un-highlighted look.
# ---------------------------------------------------------------------------
$ws.Rows("79:83").ClearFormats()

$ws.Range("A41:G41").Copy()
$ws.Range("A79:G79").PasteSpecial(-4122)

$ws.Range("A42:G42").Copy()
$ws.Range("A80:G80").PasteSpecial(-4122)

$ws.Range("A43:G43").Copy()
$ws.Range("A81:G81").PasteSpecial(-4122)

$ws.Range("A44:G44").Copy()
$ws.Range("A82:G82").PasteSpecial(-4122)

$ws.Range("A45:G45").Copy()
$ws.Range("A83:G83").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Re-apply the AutoFilter on column F with the new criterion
#    ("State Regulatory Compliance" instead of "PL Admin Tools"). This
#    naturally recomputes which rows are hidden to match the new filter.
# ---------------------------------------------------------------------------
$ws.Range("F1:F140").AutoFilter(1, @("State Regulatory Compliance"), 7)

# ---------------------------------------------------------------------------
# 3. Update the active selection/view to F34 (clearing the old selection of
#    A79:XFD83 and the scrolled topLeftCell).
# ---------------------------------------------------------------------------
$ws.Range("F34").Select()
